$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.162.33'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").Value = '2.613.88'
$ws.Range("E3").Value = '  -2.49%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.80'
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.58'
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.80%  '
$ws.Range("D9").Value = '2.612.30'
$ws.Range("E9").Value = '  -2.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.108'
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.51'
$ws.Range("E11").Value = '  -3.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.372'
$ws.Range("E12").Value = '  +4.08%  '
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.13'
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("D15").Value = '3.078.82'
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").Value = '62.998.88'
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").Value = '2.624.28'
$ws.Range("E18").Value = '  -1.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.47'
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.50'
$ws.Range("E20").Value = '  +1.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '341.40'
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.85'
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.72'
$ws.Range("E24").Value = '  -1.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.99'
$ws.Range("E25").Value = '  -2.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.69'
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.59'
$ws.Range("E27").Value = '  +2.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.00'
$ws.Range("E28").Value = '  +5.34%  '
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '543.57'
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.162'
$ws.Range("E30").Value = '  -2.51%  '
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.81'
$ws.Range("E32").Value = '  -0.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.03'
$ws.Range("E33").Value = '  +1.43%  '
$ws.Range("D34").Value = '0.0₃0838'
$ws.Range("E34").Value = '  +2.46%  '
$ws.Range("E35").Value = '  -5.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.21'
$ws.Range("E36").Value = '  +1.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '168.43'
$ws.Range("E37").Value = '  -2.70%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.402'
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.88'
$ws.Range("E41").Value = '  -1.91%  '
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '165.00'
$ws.Range("E43").Value = '  -5.59%  '
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("E45").Value = '  -0.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.71'
$ws.Range("E46").Value = '  -3.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0561'
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.622'
$ws.Range("E48").Value = '  -2.38%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.92'
$ws.Range("E51").Value = '  +10.35%  '
